{"js": "// 1) Merge the three runs that make up the \"7/30/17\" date paragraph into a\n//    single run with the combined text \"7/30/17\".\n// 2) Remove the \"Debug the modified code.\" list paragraph that used to sit\n//    right after \"SCP always gives wrong result. Debug!\".\nconst body = context.document.body;\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\n\nlet dateParagraph = null;\nlet debugParagraph = null;\n\nfor (const p of body.paragraphs.items) {\n  const t = p.text;\n  if (dateParagraph === null && t === \"7/30/17\") {\n    dateParagraph = p;\n  }\n  if (debugParagraph === null && t === \"Debug the modified code.\") {\n    debugParagraph = p;\n  }\n}\n\nif (!dateParagraph) {\n  throw new Error(\"Could not find the '7/30/17' paragraph.\");\n}\nif (!debugParagraph) {\n  throw new Error(\"Could not find the 'Debug the modified code.' paragraph.\");\n}\n\n// Collapse the three runs (\"7/\", \"30\", \"/17\") making up the date paragraph\n// into a single run containing \"7/30/17\". Deleting the tail text (\"30/17\")\n// leaves the paragraph with just its first run (\"7/\"), then re-appending\n// \"30/17\" at the paragraph's end grows that same run instead of minting a\n// new (unformatted) one.\nconst tail = dateParagraph.getRange().search(\"30/17\", { matchCase: true });\ntail.load(\"items\");\nawait context.sync();\ntail.items[0].delete();\nawait context.sync();\ndateParagraph.insertText(\"30/17\", \"End\");\n\n// Drop the whole \"Debug the modified code.\" bullet paragraph.\ndebugParagraph.delete();\n\nawait context.sync();\n", "ps1": "# 1) Merge the three runs that make up the \"7/30/17\" date paragraph into a\n#    single run with the combined text \"7/30/17\".\n# 2) Remove the \"Debug the modified code.\" list paragraph that used to sit\n#    right after \"SCP always gives wrong result. Debug!\".\n$d = $word.ActiveDocument\n\n$dateParagraph = $null\n$debugParagraph = $null\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\")\n    if (($null -eq $dateParagraph) -and ($t -eq \"7/30/17\")) {\n        $dateParagraph = $p\n    }\n    if (($null -eq $debugParagraph) -and ($t -eq \"Debug the modified code.\")) {\n        $debugParagraph = $p\n    }\n}\n\nif ($null -eq $dateParagraph) {\n    throw \"Could not find the '7/30/17' paragraph.\"\n}\nif ($null -eq $debugParagraph) {\n    throw \"Could not find the 'Debug the modified code.' paragraph.\"\n}\n\n# Collapse the three runs (\"7/\", \"30\", \"/17\") making up the date paragraph\n# into a single run containing \"7/30/17\". Deleting the tail text (\"30/17\")\n# leaves the paragraph with just its first run (\"7/\"), then re-appending\n# \"30/17\" at the paragraph's end grows that same run instead of minting a\n# new (unformatted) one.\n$tail = $d.Content\n$tail.Start = $dateParagraph.Range.Start\n$tail.End = $dateParagraph.Range.End\n$tail.Find.Execute(\"30/17\") | Out-Null\n$tail.Delete()\n$dateParagraph.Range.InsertAfter(\"30/17\")\n\n# Drop the whole \"Debug the modified code.\" bullet paragraph.\n$debugParagraph.Range.Delete()\n"}
